$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$co = $ws.ChartObjects(2)

# Move "Chart 2" down/left while keeping its current size, so the
# two-cell anchor ends up at col 0 (+257174 EMU), row 10 (+28574 EMU)
# through col 8 (+295274 EMU), row 25 (+57149 EMU).
$co.Left = 20.2499212598425
$co.Top = 144.749921259843
